$wb = $excel.ActiveWorkbook

# Duplicate the S5-N7 sheet: the copy becomes the HiSeq variant (keeps the
# original S5 index sequences), inserted before the original sheet.
$orig = $wb.Worksheets.Item("S5-N7")
$orig.Copy($orig)

# The new copy is now the first sheet in the workbook -> rename it to HiSeq
$hiseq = $wb.Worksheets.Item(1)
$hiseq.Name = "S5-N7_HiSeq"
$hiseq.Activate()
$hiseq.Range("F1").Select()

# The original sheet becomes the MiSeq variant: rename + swap in the MiSeq
# (reverse-complement) index sequences for the S50x rows (col B, rows 26-49)
$miseq = $wb.Worksheets.Item("S5-N7")
$miseq.Name = "S5-N7_MiSeq"

$miseq.Range("B26").Value = "GCGATCTA"
$miseq.Range("B27").Value = "ATAGAGAG"
$miseq.Range("B28").Value = "AGAGGATA"
$miseq.Range("B29").Value = "TCTACTCT"
$miseq.Range("B30").Value = "CTCCTTAC"
$miseq.Range("B31").Value = "TATGCAGT"
$miseq.Range("B32").Value = "TACTCCTT"
$miseq.Range("B33").Value = "AGGCTTAG"
$miseq.Range("B34").Value = "GAGTAGCC"
$miseq.Range("B35").Value = "GTCTGAGG"
$miseq.Range("B36").Value = "CGTAAGGA"
$miseq.Range("B37").Value = "CCACGCGT"
$miseq.Range("B38").Value = "GGAGTTCC"
$miseq.Range("B39").Value = "CATGGCCA"
$miseq.Range("B40").Value = "AATCTCTC"
$miseq.Range("B41").Value = "TAACCGCG"
$miseq.Range("B42").Value = "TGGCGGTC"
$miseq.Range("B43").Value = "CCATCTTA"
$miseq.Range("B44").Value = "ATGTCAAT"
$miseq.Range("B45").Value = "AGTTGGCT"
$miseq.Range("B46").Value = "ACCTAGTA"
$miseq.Range("B47").Value = "AACCGTGA"
$miseq.Range("B48").Value = "TCATTACA"
$miseq.Range("B49").Value = "CTGACGTG"

$miseq.Activate()
$miseq.Range("C25").Select()

